# Apply the crypto price/volume refresh described in the commit diff.
# Column D (Price) and column E (Volume(1h)) are stored as plain text in the
# sheet (inline strings), so number-looking price values are written with a
# leading apostrophe to force Excel to keep them as text instead of numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '63.359.89'
$ws.Range("E2").Value = '  +4.34%  '

# Row 3
$ws.Range("D3").Value = '3.489.89'
$ws.Range("E3").Value = '  +3.87%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").Value = '''585.72'
$ws.Range("E5").Value = '  +2.91%  '

# Row 6
$ws.Range("D6").Value = '''147.94'
$ws.Range("E6").Value = '  +7.74%  '

# Row 7
$ws.Range("E7").Value = '  -0.07%  '

# Row 8
$ws.Range("D8").Value = '''0.479'
$ws.Range("E8").Value = '  +1.66%  '

# Row 9
$ws.Range("D9").Value = '''7.71'
$ws.Range("E9").Value = '  +0.55%  '

# Row 10
$ws.Range("E10").Value = '  +4.72%  '

# Row 11
$ws.Range("D11").Value = '''0.400'

# Row 12
$ws.Range("D12").Value = '4.085.42'
$ws.Range("E12").Value = '  +3.82%  '

# Row 13
$ws.Range("D13").Value = '''29.87'
$ws.Range("E13").Value = '  +7.34%  '

# Row 14
$ws.Range("E14").Value = '  -0.06%  '

# Row 15
$ws.Range("D15").Value = '3.484.56'
$ws.Range("E15").Value = '  +3.68%  '

# Row 16
$ws.Range("E16").Value = '  +4.36%  '

# Row 17
$ws.Range("D17").Value = '63.342.22'
$ws.Range("E17").Value = '  +4.08%  '

# Row 18
$ws.Range("D18").Value = '''6.30'
$ws.Range("E18").Value = '  +3.44%  '

# Row 19
$ws.Range("D19").Value = '''14.40'
$ws.Range("E19").Value = '  +6.90%  '

# Row 20
$ws.Range("D20").Value = '''9.41'
$ws.Range("E20").Value = '  +6.35%  '

# Row 21
$ws.Range("D21").Value = '''393.20'
$ws.Range("E21").Value = '  +2.69%  '

# Row 22
$ws.Range("E22").Value = '  +3.48%  '

# Row 23
$ws.Range("D23").Value = '''75.16'
$ws.Range("E23").Value = '  -0.27%  '

# Row 24
$ws.Range("E24").Value = '  -0.06%  '

# Row 25
$ws.Range("D25").Value = '''0.0000120'
$ws.Range("E25").Value = '  +10.03%  '

# Row 26
$ws.Range("D26").Value = '3.630.65'
$ws.Range("E26").Value = '  +3.76%  '

# Row 27
$ws.Range("D27").Value = '''0.184'
$ws.Range("E27").Value = '  -3.22%  '

# Row 28
$ws.Range("D28").Value = '''7.85'
$ws.Range("E28").Value = '  +10.57%  '

# Row 29
$ws.Range("E29").Value = '  -0.18%  '

# Row 30
$ws.Range("D30").Value = '''8.27'
$ws.Range("E30").Value = '  +5.77%  '

# Row 31
$ws.Range("E31").Value = '  +2.47%  '

# Row 32
$ws.Range("D32").Value = '''1.43'
$ws.Range("E32").Value = '  +7.19%  '

# Row 33
$ws.Range("D33").Value = '''1.00'
$ws.Range("E33").Value = '  -0.01%  '

# Row 34
$ws.Range("D34").Value = '''23.87'
$ws.Range("E34").Value = '  +4.13%  '

# Row 35
$ws.Range("D35").Value = '''32.59'
$ws.Range("E35").Value = '  +29.94%  '

# Row 36
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '''5.35'
$ws.Range("E36").Value = '  +9.36%  '

# Row 37
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").Value = '''7.18'
$ws.Range("E37").Value = '  +5.54%  '

# Row 38
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = '''1.58'
$ws.Range("E38").Value = '  +10.21%  '

# Row 39
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").Value = '''171.38'
$ws.Range("E39").Value = '  +2.25%  '

# Row 40
$ws.Range("D40").Value = '3.527.29'
$ws.Range("E40").Value = '  +3.90%  '

# Row 41
$ws.Range("D41").Value = '''0.0769'
$ws.Range("E41").Value = '  +2.10%  '

# Row 42
$ws.Range("D42").Value = '''0.804'
$ws.Range("E42").Value = '  +4.78%  '

# Row 43
$ws.Range("D43").Value = '''1.74'
$ws.Range("E43").Value = '  +7.95%  '

# Row 44
$ws.Range("D44").Value = '''4.51'
$ws.Range("E44").Value = '  +4.30%  '

# Row 45
$ws.Range("D45").Value = '''42.47'
$ws.Range("E45").Value = '  +0.60%  '

# Row 46
$ws.Range("E46").Value = '  +10.79%  '

# Row 47
$ws.Range("D47").Value = '2.618.97'
$ws.Range("E47").Value = '  +7.30%  '

# Row 48
$ws.Range("D48").Value = '''23.91'
$ws.Range("E48").Value = '  +8.33%  '

# Row 49
$ws.Range("E49").Value = '  +18.86%  '

# Row 50
$ws.Range("E50").Value = '  +2.53%  '

# Row 51
$ws.Range("E51").Value = '  +5.87%  '
